$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Date line
Replace-Text "2024-09-04 Wednesday" "2024-09-05 Thursday"

# Table 1 (the only table in the document)
$tbl = $d.Tables.Item(1)

# Row 1
$tbl.Cell(1,1).Range.Text = "30÷6="
$tbl.Cell(1,2).Range.Text = "40÷3="
$tbl.Cell(1,3).Range.Text = "22÷9="
$tbl.Cell(1,4).Range.Text = "67÷6="
$tbl.Cell(1,5).Range.Text = "92÷2="

# Row 5
$tbl.Cell(5,1).Range.Text = "86÷7="
$tbl.Cell(5,2).Range.Text = "15÷2="
$tbl.Cell(5,3).Range.Text = "81÷4="
$tbl.Cell(5,4).Range.Text = "95÷7="
$tbl.Cell(5,5).Range.Text = "32÷8="

# Row 9
$tbl.Cell(9,1).Range.Text = "61÷9="
$tbl.Cell(9,2).Range.Text = "50÷9="
$tbl.Cell(9,3).Range.Text = "56÷3="
$tbl.Cell(9,4).Range.Text = "30÷6="
$tbl.Cell(9,5).Range.Text = "84÷6="

# Row 13
$tbl.Cell(13,1).Range.Text = "57÷8="
$tbl.Cell(13,2).Range.Text = "24÷6="
$tbl.Cell(13,3).Range.Text = "67÷8="
$tbl.Cell(13,4).Range.Text = "57÷3="
$tbl.Cell(13,5).Range.Text = "45÷9="

# Row 17
$tbl.Cell(17,1).Range.Text = "94÷4="
$tbl.Cell(17,2).Range.Text = "94÷2="
$tbl.Cell(17,3).Range.Text = "65÷6="
$tbl.Cell(17,4).Range.Text = "66÷2="
$tbl.Cell(17,5).Range.Text = "74÷6="
